# Fix Training Data Issue (#48)
# Data was taken from 1 day off due to way NBA stats were shown.
# The "Date" column (BF) on Sheet1 held the literal text "5-2-2007-08"
# for every data row (rows 2-31). Correct it to the literal text
# "2008-05-02".
#
# Note: assigning a date-looking string straight to Range.Value makes
# this host auto-convert it into a real date serial (like Excel's
# normal typed-input inference), which is not what we want here - the
# source data stores the date as plain text. Temporarily forcing the
# cells to a text number format before the assignment keeps the value
# a literal string, and ClearFormats() afterwards drops that temporary
# formatting again so the cells end up with their original (default)
# style, untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$firstRow = 2
$lastRow = 31
$col = "BF"
$oldDate = "5-2-2007-08"
$newDate = "2008-05-02"

$rng = $ws.Range($col + $firstRow + ":" + $col + $lastRow)
$rng.NumberFormat = "@"

for ($r = $firstRow; $r -le $lastRow; $r++) {
    $cell = $ws.Range($col + $r)
    if ($cell.Text -eq $oldDate) {
        $cell.Value = $newDate
    }
}

$rng.ClearFormats()
